# Add a new 'Plates' column (column E) to sheets ED3A and ED3D,
# listing the plate barcodes associated with each comparison row.
$wb = $excel.ActiveWorkbook

# --- Sheet ED3A ---
$ws1 = $wb.Worksheets.Item("ED3A")
$valuesED3A = @{
    2 = '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'']'
    3 = '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00117020'', ''BR00117021'']'
    4 = '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00118050'', ''BR00117006'']'
    5 = '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00117020'', ''BR00117021'']'
    6 = '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00118050'', ''BR00117006'']'
    7 = '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'']'
    8 = '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']'
    9 = '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'']'
    10 = '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']'
    11 = '[''BR00117020'', ''BR00117021'', ''BR00118050'', ''BR00117006'']'
    12 = '[''BR00117020'', ''BR00117021'', ''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'']'
    13 = '[''BR00117020'', ''BR00117021'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']'
    14 = '[''BR00118050'', ''BR00117006'', ''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'']'
    15 = '[''BR00118050'', ''BR00117006'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']'
    16 = '[''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']'
    17 = '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00118049'']'
    18 = '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00118049'']'
    19 = '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00117001'', ''BR00117002'']'
    20 = '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00117001'', ''BR00117002'']'
    21 = '[''BR00117020'', ''BR00117021'', ''BR00118049'']'
    22 = '[''BR00118050'', ''BR00117006'', ''BR00118049'']'
    23 = '[''BR00117020'', ''BR00117021'', ''BR00117001'', ''BR00117002'']'
    24 = '[''BR00118050'', ''BR00117006'', ''BR00117001'', ''BR00117002'']'
    25 = '[''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'', ''BR00118049'']'
    26 = '[''BR00118050'', ''BR00117006'', ''BR00117001'', ''BR00117002'']'
    27 = '[''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'', ''BR00117001'', ''BR00117002'']'
    28 = '[''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'', ''BR00117001'', ''BR00117002'']'
    29 = '[''BR00118049'', ''BR00117001'', ''BR00117002'']'
}
foreach ($row in $valuesED3A.Keys) {
    $ws1.Cells.Item($row, 5).Value = $valuesED3A[$row]
}
$ws1.Cells.Item(1, 5).Value = "Plates"
$ws1.Range("D1").Copy() | Out-Null
$ws1.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Sheet ED3D ---
$ws4 = $wb.Worksheets.Item("ED3D")
$valuesED3D = @{
    2 = '[''BR00121434'', ''BR00121435'', ''BR00121440'', ''BR00121431'', ''BR00121432'', ''BR00121433'']'
    3 = '[''BR00121434'', ''BR00121435'', ''BR00121440'', ''BR00117022'', ''BR00117023'']'
    4 = '[''BR00121434'', ''BR00121435'', ''BR00121440'', ''BR00118039'', ''BR00118040'']'
    5 = '[''BR00121431'', ''BR00121432'', ''BR00121433'', ''BR00117022'', ''BR00117023'']'
    6 = '[''BR00121431'', ''BR00121432'', ''BR00121433'', ''BR00118039'', ''BR00118040'']'
    7 = '[''BR00117022'', ''BR00117023'', ''BR00118039'', ''BR00118040'']'
}
foreach ($row in $valuesED3D.Keys) {
    $ws4.Cells.Item($row, 5).Value = $valuesED3D[$row]
}
$ws4.Cells.Item(1, 5).Value = "Plates"
$ws4.Range("D1").Copy() | Out-Null
$ws4.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Output "done"
